$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.482.26"
$ws.Range("E2").Value = "  +0.87%  "

$ws.Range("D3").Value = "1.873.22"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7188"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "239.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.60%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07796"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.38%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3072"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08242"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.82%  "

$ws.Range("D12").Value = "1.884.55"
$ws.Range("E12").Value = "  +0.15%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.238"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7215"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.60%  "

$ws.Range("D16").Value = "29.526.17"
$ws.Range("E16").Value = "  +1.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.844"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.86%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007851"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.63"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.40%  "

$ws.Range("D21").Value = "2.134.16"
$ws.Range("E21").Value = "  +1.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.748"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1566"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.38%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.012"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.936"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.355"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.96%  "

$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.337"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.090"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05255"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.200"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7185"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.001"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.677"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01868"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.42%  "

$ws.Range("E40").Value = "  -0.43%  "

$ws.Range("D41").Value = "1.180.23"
$ws.Range("E41").Value = "  +3.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9071"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.03%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.20%  "

$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4314"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5360"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.767"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.172"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.026"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.01%  "
